$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text so Excel
# doesn't silently convert the inline string into a number (diff keeps
# these as t="inlineStr"). Apply a temporary Text format, assign the
# value, then reset the style back to Normal so no stray style index is
# left behind on the cell.
$numericLooking = @(
    'D4'
    'D5'
    'D6'
    'D8'
    'D9'
    'D10'
    'D12'
    'D14'
    'D15'
    'D16'
    'D17'
    'D18'
    'D19'
    'D20'
    'D21'
    'D23'
    'D24'
    'D27'
    'D30'
    'D31'
    'D32'
    'D33'
    'D34'
    'D35'
    'D36'
    'D37'
    'D39'
    'D40'
    'D41'
    'D42'
    'D43'
    'D44'
    'D46'
    'D47'
    'D48'
    'D50'
    'D51'
)
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D4').Value = '0.9998'
$ws.Range('D5').Value = '333.48'
$ws.Range('D6').Value = '0.9996'
$ws.Range('D8').Value = '0.4097'
$ws.Range('D9').Value = '47.70'
$ws.Range('D10').Value = '0.08014'
$ws.Range('D12').Value = '21.85'
$ws.Range('D14').Value = '5.952'
$ws.Range('D15').Value = '7.098'
$ws.Range('D16').Value = '89.20'
$ws.Range('D17').Value = '1.001'
$ws.Range('D18').Value = '0.00001031'
$ws.Range('D19').Value = '0.06572'
$ws.Range('D20').Value = '17.54'
$ws.Range('D21').Value = '0.9994'
$ws.Range('D23').Value = '5.448'
$ws.Range('D24').Value = '11.31'
$ws.Range('D27').Value = '157.27'
$ws.Range('D30').Value = '5.426'
$ws.Range('D31').Value = '118.55'
$ws.Range('D32').Value = '0.9907'
$ws.Range('D33').Value = '0.09421'
$ws.Range('D34').Value = '1.431'
$ws.Range('D35').Value = '3.595'
$ws.Range('D36').Value = '5.322'
$ws.Range('D37').Value = '0.06107'
$ws.Range('D39').Value = '8.389'
$ws.Range('D40').Value = '1.181'
$ws.Range('D41').Value = '0.5823'
$ws.Range('D42').Value = '0.9992'
$ws.Range('D43').Value = '10.20'
$ws.Range('D44').Value = '0.1827'
$ws.Range('D46').Value = '2.349'
$ws.Range('D47').Value = '12.11'
$ws.Range('D48').Value = '0.5516'
$ws.Range('D50').Value = '0.07052'
$ws.Range('D51').Value = '47.83'

foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cells (words, URLs, multi-dot numbers, percent strings) are
# already safe from numeric auto-conversion.
$ws.Range('D2').Value = '29.090.36'
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('D3').Value = '1.910.13'
$ws.Range('E3').Value = '  +2.14%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  -1.21%  '
$ws.Range('E8').Value = '  +3.24%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = '1.899.71'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('E15').Value = '  -2.27%  '
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '29.088.43'
$ws.Range('E22').Value = '  +2.40%  '
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('D26').Value = '2.127.87'
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('E27').Value = '  -2.15%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('E34').Value = '  +4.02%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('E40').Value = '  +0.37%  '
$ws.Range('E41').Value = '  -2.26%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('E44').Value = '  -2.85%  '
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('E46').Value = '  +13.93%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('E49').Value = '  -1.78%  '
$ws.Range('E50').Value = '  +1.81%  '
$ws.Range('E51').Value = '  +22.70%  '
